$wb = $excel.ActiveWorkbook

# IPS Clutches & Brakes Division
$ws = $wb.Worksheets.Item("IPS Clutches & Brakes Division")
$ws.Range("D2").Value = 0.0502
$ws.Range("D3").Value = 0.0502
$ws.Range("D4").Value = 0.0502
$ws.Range("H4").Value = 0.0069
$ws.Range("I4").Value = 0.0196
$ws.Range("M4").Value = 0.0227
$ws.Range("N4").Value = 0.008
$ws.Range("O4").Value = 0.007175
$ws.Range("P4").Value = 0.007175
$ws.Range("Q4").Value = 0.021525
$ws.Range("R4").Value = 0.007175
$ws.Range("S4").Value = 0.007175
$ws.Range("T4").Value = 0.007175
$ws.Range("U4").Value = 0.021525
$ws.Range("V4").Value = 0.0861
$ws.Range("D5").Value = 0.511627906976744
$ws.Range("D6").Value = 0.511627906976744
$ws.Range("D7").Value = 0.511627906976744
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0.511627906976744
$ws.Range("P7").Value = 0.511627906976744
$ws.Range("Q7").Value = 0.511627906976744
$ws.Range("R7").Value = 0.511627906976744
$ws.Range("S7").Value = 0.511627906976744
$ws.Range("T7").Value = 0.511627906976744
$ws.Range("U7").Value = 0.511627906976744
$ws.Range("V7").Value = 0.511627906976744

# IPS Couplings Division
$ws = $wb.Worksheets.Item("IPS Couplings Division")
$ws.Range("D2").Value = 0.0541
$ws.Range("D3").Value = 0.0541
$ws.Range("D4").Value = 0.0541
$ws.Range("L4").Value = 0.0085
$ws.Range("N4").Value = 0.0049
$ws.Range("O4").Value = 0.007725
$ws.Range("P4").Value = 0.007725
$ws.Range("Q4").Value = 0.023175
$ws.Range("R4").Value = 0.007725
$ws.Range("S4").Value = 0.007725
$ws.Range("T4").Value = 0.007725
$ws.Range("U4").Value = 0.023175
$ws.Range("V4").Value = 0.0927
$ws.Range("D5").Value = 0.698412698412698
$ws.Range("D6").Value = 0.698412698412698
$ws.Range("D7").Value = 0.698412698412698
$ws.Range("H7").Value = 0.4
$ws.Range("I7").Value = 0.2143
$ws.Range("L7").Value = 0.5714
$ws.Range("M7").Value = 0.4615
$ws.Range("N7").Value = 0.9722
$ws.Range("O7").Value = 0.698412698412698
$ws.Range("P7").Value = 0.698412698412698
$ws.Range("Q7").Value = 0.698412698412698
$ws.Range("R7").Value = 0.698412698412698
$ws.Range("S7").Value = 0.698412698412698
$ws.Range("T7").Value = 0.698412698412698
$ws.Range("U7").Value = 0.698412698412698
$ws.Range("V7").Value = 0.698412698412698

# IPS Gearing Division
$ws = $wb.Worksheets.Item("IPS Gearing Division")
$ws.Range("D2").Value = 0.0786
$ws.Range("D3").Value = 0.0786
$ws.Range("D4").Value = 0.0786
$ws.Range("F4").Value = 0.0097
$ws.Range("H4").Value = 0.0114
$ws.Range("I4").Value = 0.0211
$ws.Range("K4").Value = 0.0165
$ws.Range("L4").Value = 0.0166
$ws.Range("M4").Value = 0.0445
$ws.Range("N4").Value = 0.0133
$ws.Range("O4").Value = 0.011225
$ws.Range("P4").Value = 0.011225
$ws.Range("Q4").Value = 0.033675
$ws.Range("R4").Value = 0.011225
$ws.Range("S4").Value = 0.011225
$ws.Range("T4").Value = 0.011225
$ws.Range("U4").Value = 0.033675
$ws.Range("V4").Value = 0.1347
$ws.Range("D5").Value = 0.466666666666667
$ws.Range("D6").Value = 0.466666666666667
$ws.Range("D7").Value = 0.466666666666667
$ws.Range("F7").Value = 0.5714
$ws.Range("I7").Value = 0.6364
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.4286
$ws.Range("N7").Value = 0.2
$ws.Range("O7").Value = 0.466666666666667
$ws.Range("P7").Value = 0.466666666666667
$ws.Range("Q7").Value = 0.466666666666667
$ws.Range("R7").Value = 0.466666666666667
$ws.Range("S7").Value = 0.466666666666667
$ws.Range("T7").Value = 0.466666666666667
$ws.Range("U7").Value = 0.466666666666667
$ws.Range("V7").Value = 0.466666666666667

# IPS Industrial Components Divi
$ws = $wb.Worksheets.Item("IPS Industrial Components Divi")
$ws.Range("D2").Value = 0.0632
$ws.Range("D3").Value = 0.0632
$ws.Range("D4").Value = 0.0632
$ws.Range("I4").Value = 0.0217
$ws.Range("N4").Value = 0.0088
$ws.Range("O4").Value = 0.009025
$ws.Range("P4").Value = 0.009025
$ws.Range("Q4").Value = 0.027075
$ws.Range("R4").Value = 0.009025
$ws.Range("S4").Value = 0.009025
$ws.Range("T4").Value = 0.009025
$ws.Range("U4").Value = 0.027075
$ws.Range("V4").Value = 0.1083
$ws.Range("D5").Value = 0.5625
$ws.Range("D6").Value = 0.5625
$ws.Range("D7").Value = 0.5625
$ws.Range("N7").Value = 0.7143
$ws.Range("O7").Value = 0.5625
$ws.Range("P7").Value = 0.5625
$ws.Range("Q7").Value = 0.5625
$ws.Range("R7").Value = 0.5625
$ws.Range("S7").Value = 0.5625
$ws.Range("T7").Value = 0.5625
$ws.Range("U7").Value = 0.5625
$ws.Range("V7").Value = 0.5625

# IPS Segment Functions
$ws = $wb.Worksheets.Item("IPS Segment Functions")
$ws.Range("D2").Value = 0.059
$ws.Range("D3").Value = 0.059
$ws.Range("D4").Value = 0.059
$ws.Range("J4").Value = 0.0067
$ws.Range("M4").Value = 0.0244
$ws.Range("N4").Value = 0.0111
$ws.Range("O4").Value = 0.008425
$ws.Range("P4").Value = 0.008425
$ws.Range("Q4").Value = 0.025275
$ws.Range("R4").Value = 0.008425
$ws.Range("S4").Value = 0.008425
$ws.Range("T4").Value = 0.008425
$ws.Range("U4").Value = 0.025275
$ws.Range("V4").Value = 0.1011
$ws.Range("D5").Value = 0.654545454545455
$ws.Range("D6").Value = 0.654545454545455
$ws.Range("D7").Value = 0.654545454545455
$ws.Range("N7").Value = 0.75
$ws.Range("O7").Value = 0.654545454545455
$ws.Range("P7").Value = 0.654545454545455
$ws.Range("Q7").Value = 0.654545454545455
$ws.Range("R7").Value = 0.654545454545455
$ws.Range("S7").Value = 0.654545454545455
$ws.Range("T7").Value = 0.654545454545455
$ws.Range("U7").Value = 0.654545454545455
$ws.Range("V7").Value = 0.654545454545455

# L1_IPS
$ws = $wb.Worksheets.Item("L1_IPS")
$ws.Range("D2").Value = 0.0592
$ws.Range("D3").Value = 0.0592
$ws.Range("D4").Value = 0.0592
$ws.Range("K4").Value = 0.0102
$ws.Range("N4").Value = 0.0091
$ws.Range("O4").Value = 0.00845833333333333
$ws.Range("P4").Value = 0.00845833333333333
$ws.Range("Q4").Value = 0.025375
$ws.Range("R4").Value = 0.00845833333333333
$ws.Range("S4").Value = 0.00845833333333333
$ws.Range("T4").Value = 0.00845833333333333
$ws.Range("U4").Value = 0.025375
$ws.Range("V4").Value = 0.1015
$ws.Range("D5").Value = 0.600896860986547
$ws.Range("D6").Value = 0.600896860986547
$ws.Range("D7").Value = 0.600896860986547
$ws.Range("F7").Value = 0.5526
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.5435
$ws.Range("N7").Value = 0.8448
$ws.Range("O7").Value = 0.600896860986547
$ws.Range("P7").Value = 0.600896860986547
$ws.Range("Q7").Value = 0.600896860986547
$ws.Range("R7").Value = 0.600896860986547
$ws.Range("S7").Value = 0.600896860986547
$ws.Range("T7").Value = 0.600896860986547
$ws.Range("U7").Value = 0.600896860986547
$ws.Range("V7").Value = 0.600896860986547

Write-Output "Applied CVD updates"